$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46 - "Invert a binary tree"
$ws.Range("C46").Value = "Recursive solution, swapping left and right subtrees on the way up"
$ws.Range("D46").Value = "Linear, visit each node once"
$ws.Range("E46").Value = "Call stack"
$ws.Range("F46").Value = "8/14/2025"
$ws.Range("G46").Value = "10 minutes"
$ws.Range("H46").Value = "10 minutes"
$ws.Range("I46").Value = "N/A"
$ws.Range("J46").Value = "Yes"
$ws.Range("K46").Value = "Started with subtree, and kept reversing node till solution"
$ws.Range("L46").Value = "No"
$ws.Range("M46").Value = "No"
$ws.Range("N46").Value = "Almost the same solution as book"
$ws.Range("O46").Value = "No"
$ws.Range("P46").Value = 4
$ws.Range("Q46").Value = 4
$ws.Range("R46").Value = 4
$ws.Range("S46").Value = 4
$ws.Rows.Item(46).RowHeight = 60

# Row 47 - "Evaluate expression tree"
$ws.Range("C47").Value = "Bottom-up, evaluate based on kind of node"
$ws.Range("D47").Value = "Linear, number of nodes"
$ws.Range("E47").Value = "Call stack height"
$ws.Range("F47").Value = "8/14/2025"
$ws.Range("G47").Value = "10 minutes"
$ws.Range("H47").Value = "10 minutes"
$ws.Range("I47").Value = "No"
$ws.Range("J47").Value = "Yes"
$ws.Range("K47").Value = "All good"
$ws.Range("L47").Value = "No"
$ws.Range("M47").Value = "No"
$ws.Range("N47").Value = "Almost the same solution as book"
$ws.Range("O47").Value = "Either memorise math.prod or keep a helper function handy"
$ws.Range("P47").Value = 4
$ws.Range("Q47").Value = 4
$ws.Range("R47").Value = 4
$ws.Range("S47").Value = 4
$ws.Rows.Item(47).RowHeight = 60

$ws.Range("T47").Select()
